# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# -> Update column G ("K") values for rows 2-15 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 2
    6  = 2
    7  = 2
    8  = 4
    9  = 2
    10 = 3
    11 = 1
    12 = 3
    13 = 3
    14 = 1
    15 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
